# Adds newly scraped rows (30/12/2025, scrape run 19:40) to the three
# "horarios 141" sheets, and refreshes the "Última actualización" /
# "Total filas" header cells (A2 / A3) on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912   (columns: A -, B Hora_Scrap, C Hora_Llegada, D Linea,
#                     E Minutos, F Parada, G Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 19:40:33"
$ws1.Range("A3").Value = "Total filas: 581"

$rows1 = @(
    ,@(566, "19:40:22", "19:50", "11X44_ETCHEVERRY",          10, "LP1912", "30/12/2025")
    ,@(567, "19:40:22", "19:50", "16_P MOR-SANTA ANA",        10, "LP1912", "30/12/2025")
    ,@(568, "19:40:22", "19:51", "81_EL PELIGRO",             11, "LP1912", "30/12/2025")
    ,@(569, "19:40:22", "19:59", "17_ROMERO",                 19, "LP1912", "30/12/2025")
    ,@(570, "19:40:22", "20:00", "14_ABASTO",                 20, "LP1912", "30/12/2025")
    ,@(571, "19:40:22", "20:01", "16_SANTA ANA",               21, "LP1912", "30/12/2025")
    ,@(572, "19:40:22", "20:07", "10_OLMOS",                  27, "LP1912", "30/12/2025")
    ,@(573, "19:40:22", "20:09", "15_ABASTO",                 29, "LP1912", "30/12/2025")
    ,@(574, "19:40:22", "20:10", "16_P MOR-167 Y 521",        30, "LP1912", "30/12/2025")
    ,@(575, "19:40:22", "20:12", "23_HERNANDEZ",              32, "LP1912", "30/12/2025")
    ,@(576, "19:40:22", "20:21", "26_HERNANDEZ",              41, "LP1912", "30/12/2025")
    ,@(577, "19:40:22", "20:22", "11_ETCHEVERRY",             42, "LP1912", "30/12/2025")
    ,@(578, "19:40:22", "20:23", "215A_EL PATO",              43, "LP1912", "30/12/2025")
    ,@(579, "19:40:22", "20:52", "15_ABASTO",                 72, "LP1912", "30/12/2025")
    ,@(580, "19:40:22", "20:57", "23_HERNANDEZ",              77, "LP1912", "30/12/2025")
    ,@(581, "19:40:22", "21:04", "84_COLONIA URQUIZA-ESC 49", 84, "LP1912", "30/12/2025")
    ,@(582, "19:40:22", "21:07", "215B_EL PATO",               87, "LP1912", "30/12/2025")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215   (columns: A -, B Fecha, C Hora_Scrap, D Hora_Llegada,
#                         E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 19:40:33"
$ws2.Range("A3").Value = "Total filas: 39"

$rows2 = @(
    ,@(39, "30/12/2025", "19:40:22", "20:23", "215A_EL PATO", 43, "LP1912")
    ,@(40, "30/12/2025", "19:40:22", "21:07", "215B_EL PATO", 87, "LP1912")
)

foreach ($row in $rows2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173   (columns: A -, B Fecha, C Hora_Scrap, D Hora_Llegada,
#                        E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 19:40:33"
$ws3.Range("A3").Value = "Total filas: 72"

$rows3 = @(
    ,@(73, "30/12/2025", "19:40:28", "19:54", "215C_LA PLATA", 14, "L6203")
)

foreach ($row in $rows3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}
